$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, 'Biltin', 1, 'Biltin 20mg Tablet 30''s', '30''s')
    ,@(4, 'Desodin', 2, 'Desodin 60ml Syrup', '60 ml')
    ,@(5, 'Dinafex', 3, 'Dinafex 180mg Tablet', '30''s')
    ,@(5, 'Dinafex', 4, 'Dinafex 120mg Tablet', '30''s')
    ,@(5, 'Dinafex', 5, 'Dinafex 60mg Tablet', '30''s')
    ,@(6, 'Dorenta', 6, 'Dorenta 50mg Tablet', '50''s')
    ,@(7, 'Etorix', 7, 'Etorix 90mg Tablet', '30''s')
    ,@(7, 'Etorix', 8, 'Etorix 120mg Tablet', '20''s')
    ,@(7, 'Etorix', 9, 'Etorix 60mg Tablet - 40''s', '40''s')
    ,@(8, 'Fenobac', 10, 'Fenobac 100ml Syrup', '100ml')
    ,@(9, 'Flucloxin', 11, 'Flucloxin 500mg Capsule', '30 ''s')
    ,@(9, 'Flucloxin', 12, 'Flucloxin 500mg Capsule - 36''s', '36 ''s')
    ,@(10, 'Geminox', 13, 'Geminox 320mg Tablet - 8''s', '8 ''s')
    ,@(11, 'Ketonic', 14, 'Ketonic 30mg Injection', '5 ''s')
    ,@(11, 'Ketonic', 15, 'Ketonic 10mg Tablet', '20''s')
    ,@(11, 'Ketonic', 16, 'Ketonic 30mg IM/IV Injection - 4''s', '4''s')
    ,@(12, 'Kynol', 17, 'Kynol TR 200mg Capsule', '30 ''s')
    ,@(12, 'Kynol', 18, 'Kynol D 25mg Tablet', '60 ''s')
    ,@(12, 'Kynol', 19, 'Kynol TR 100mg Capsule', '50 ''s')
    ,@(17, 'Naprox', 20, 'Naprox Plus 500mg Tablet - 30''s', '30 ''s')
    ,@(19, 'Oradin', 21, 'Oradin Plus Tablet - 40''s', '40 ''s')
    ,@(20, 'Osticare', 22, 'Osticare Tablet 24''s', '24''s')
    ,@(23, 'Rupaday', 23, 'Rupaday Oral Solution 60ml', '1''s')
    ,@(24, 'Sk-Mox', 24, 'Sk-Mox 500mg Capsule', '48 ''s')
    ,@(35, 'Zithrox', 25, 'Zithrox 15ml Suspension', '15 ml')
    ,@(35, 'Zithrox', 26, 'Zithrox 500mg Tablet', '6 ''s')
    ,@(35, 'Zithrox', 27, 'Zithrox 30ml Dry Suspension', '30ml')
    ,@(35, 'Zithrox', 28, 'Zithrox 250mg Tablet - 6''s', '6''s')
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row = $row + 1
}
